# "Yearly coverage in scenario 1"
#
# On the MarketShare sheet, the yearly coverage values (1) that were
# entered on row 3 (Old Product B (SOC)) in columns L:Z (years 2026-2040)
# actually belonged on row 2 (New Product A) instead. Row 3 keeps its
# existing D:K values (years 2018-2025).
#
# Also re-point the active sheet/selection at MarketShare (previously the
# "Platform Coverage" sheet was the active/selected tab).

$wb  = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("MarketShare")

# Fill in the yearly coverage for "New Product A" (row 2), years 2026-2040.
$ws2.Range("L2:Z2").Value = 1

# Remove the same values that were mistakenly placed on row 3.
$ws2.Range("L3:Z3").ClearContents()

# MarketShare becomes the active sheet/tab, with Z2 selected.
$ws2.Activate()
$ws2.Range("Z2").Select()
